# Fruta / hortaliza, semanal
# Weekly refresh of the "Fruta, Feria Lagunitas de Puerto Montt - Pera" price
# series: existing rows 276-299 get their date / volume / price figures
# (and a few variety/quality labels) corrected, and two brand-new
# observations are inserted ahead of the former last row (old row 300),
# which is pushed down to row 302. The sheet's used range grows from
# A1:T300 to A1:T302 automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 276-299 with corrected values ---
# Row 276
$ws.Range("D276").Value = 44769
$ws.Range("M276").Value = 100
$ws.Range("N276").Value = 16000
$ws.Range("O276").Value = 16000
$ws.Range("P276").Value = 16000
$ws.Range("S276").Value = 1067

# Row 277
$ws.Range("D277").Value = 44769
$ws.Range("L277").Value = 'Segunda'
$ws.Range("M277").Value = 100
$ws.Range("N277").Value = 13000
$ws.Range("P277").Value = 13000
$ws.Range("S277").Value = 867

# Row 278
$ws.Range("D278").Value = 44596
$ws.Range("N278").Value = 14000
$ws.Range("O278").Value = 15000
$ws.Range("P278").Value = 14500
$ws.Range("S278").Value = 967

# Row 279
$ws.Range("D279").Value = 44340
$ws.Range("L279").Value = 'Primera'
$ws.Range("O279").Value = 13000
$ws.Range("P279").Value = 12500
$ws.Range("S279").Value = 833

# Row 280
$ws.Range("D280").Value = 44425
$ws.Range("M280").Value = 400
$ws.Range("N280").Value = 17000
$ws.Range("O280").Value = 17000
$ws.Range("P280").Value = 17000
$ws.Range("S280").Value = 1133

# Row 281
$ws.Range("D281").Value = 44425
$ws.Range("M281").Value = 200
$ws.Range("N281").Value = 12000
$ws.Range("O281").Value = 12000
$ws.Range("P281").Value = 12000
$ws.Range("S281").Value = 800

# Row 282
$ws.Range("D282").Value = 44512
$ws.Range("M282").Value = 600
$ws.Range("N282").Value = 15000
$ws.Range("O282").Value = 16000
$ws.Range("P282").Value = 15500
$ws.Range("S282").Value = 1033

# Row 283
$ws.Range("D283").Value = 44512
$ws.Range("L283").Value = 'Segunda'
$ws.Range("N283").Value = 14000
$ws.Range("O283").Value = 14000
$ws.Range("P283").Value = 14000
$ws.Range("S283").Value = 933

# Row 284
$ws.Range("D284").Value = 44397
$ws.Range("M284").Value = 400
$ws.Range("N284").Value = 12000
$ws.Range("O284").Value = 12000
$ws.Range("P284").Value = 12000
$ws.Range("S284").Value = 800

# Row 285
$ws.Range("D285").Value = 44497
$ws.Range("K285").Value = 'Packham''s Triumph'
$ws.Range("M285").Value = 300
$ws.Range("N285").Value = 15000
$ws.Range("O285").Value = 16000
$ws.Range("P285").Value = 15500
$ws.Range("S285").Value = 1033

# Row 286
$ws.Range("D286").Value = 44285
$ws.Range("M286").Value = 300
$ws.Range("N286").Value = 14000
$ws.Range("O286").Value = 15000
$ws.Range("P286").Value = 14500
$ws.Range("S286").Value = 967

# Row 287
$ws.Range("D287").Value = 44362
$ws.Range("K287").Value = 'Forelle'
$ws.Range("N287").Value = 12000
$ws.Range("O287").Value = 13000
$ws.Range("P287").Value = 12500
$ws.Range("S287").Value = 833

# Row 288
$ws.Range("D288").Value = 44362
$ws.Range("L288").Value = 'Primera'
$ws.Range("M288").Value = 500
$ws.Range("N288").Value = 13500
$ws.Range("O288").Value = 14000
$ws.Range("P288").Value = 13750
$ws.Range("S288").Value = 917

# Row 289
$ws.Range("D289").Value = 44747
$ws.Range("M289").Value = 400
$ws.Range("N289").Value = 15000
$ws.Range("O289").Value = 15000
$ws.Range("P289").Value = 15000
$ws.Range("S289").Value = 1000

# Row 290
$ws.Range("D290").Value = 44747
$ws.Range("L290").Value = 'Segunda'
$ws.Range("M290").Value = 200
$ws.Range("N290").Value = 12000
$ws.Range("O290").Value = 12000
$ws.Range("P290").Value = 12000
$ws.Range("S290").Value = 800

# Row 291
$ws.Range("D291").Value = 44357
$ws.Range("M291").Value = 120
$ws.Range("N291").Value = 13500
$ws.Range("O291").Value = 14000
$ws.Range("P291").Value = 13750
$ws.Range("S291").Value = 917

# Row 292
$ws.Range("D292").Value = 44551
$ws.Range("L292").Value = 'Primera'
$ws.Range("M292").Value = 500
$ws.Range("N292").Value = 13000
$ws.Range("O292").Value = 14000
$ws.Range("P292").Value = 13500
$ws.Range("S292").Value = 900

# Row 293
$ws.Range("D293").Value = 44757
$ws.Range("M293").Value = 300
$ws.Range("N293").Value = 16000
$ws.Range("O293").Value = 16000
$ws.Range("P293").Value = 16000
$ws.Range("S293").Value = 1067

# Row 294
$ws.Range("D294").Value = 44757
$ws.Range("L294").Value = 'Segunda'
$ws.Range("M294").Value = 300
$ws.Range("N294").Value = 12000
$ws.Range("O294").Value = 12000
$ws.Range("P294").Value = 12000
$ws.Range("S294").Value = 800

# Row 295
$ws.Range("D295").Value = 44547
$ws.Range("K295").Value = 'Packham''s Triumph'
$ws.Range("M295").Value = 600
$ws.Range("N295").Value = 13000
$ws.Range("O295").Value = 14000
$ws.Range("P295").Value = 13500
$ws.Range("S295").Value = 900

# Row 296
$ws.Range("D296").Value = 44438
$ws.Range("M296").Value = 200
$ws.Range("N296").Value = 16000
$ws.Range("O296").Value = 16000
$ws.Range("P296").Value = 16000
$ws.Range("S296").Value = 1067

# Row 297
$ws.Range("D297").Value = 44355
$ws.Range("K297").Value = 'Forelle'
$ws.Range("M297").Value = 200
$ws.Range("O297").Value = 13000
$ws.Range("P297").Value = 12500
$ws.Range("S297").Value = 833

# Row 298
$ws.Range("D298").Value = 44355
$ws.Range("K298").Value = 'Packham''s Triumph'
$ws.Range("M298").Value = 300
$ws.Range("N298").Value = 13500
$ws.Range("O298").Value = 14000
$ws.Range("P298").Value = 13750
$ws.Range("S298").Value = 917

# Row 299
$ws.Range("D299").Value = 44391
$ws.Range("M299").Value = 300
$ws.Range("N299").Value = 12000
$ws.Range("O299").Value = 12000
$ws.Range("P299").Value = 12000
$ws.Range("S299").Value = 800

# --- Insert two brand-new observation rows ahead of the old last row ---
# (old row 300 gets pushed down to row 302, unchanged)
$ws.Rows.Item(300).Insert()
$ws.Rows.Item(301).Insert()

# New row 300
$ws.Range("A300").Value = 4
$ws.Range("B300").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C300").Value = 'Los Lagos'
$ws.Range("D300").Value = 44358
$ws.Range("E300").Value = 10
$ws.Range("F300").Value = 'Fruta'
$ws.Range("G300").Value = 100104
$ws.Range("H300").Value = 'Frutos de pepita'
$ws.Range("I300").Value = 100104005
$ws.Range("J300").Value = 'Pera'
$ws.Range("K300").Value = 'Forelle'
$ws.Range("L300").Value = 'Primera'
$ws.Range("M300").Value = 400
$ws.Range("N300").Value = 12000
$ws.Range("O300").Value = 13000
$ws.Range("P300").Value = 12500
$ws.Range("Q300").Value = '$/caja 15 kilos empedrada'
$ws.Range("R300").Value = "Región de O'Higgins"
$ws.Range("S300").Value = 833
$ws.Range("T300").Value = 15

# New row 301
$ws.Range("A301").Value = 4
$ws.Range("B301").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C301").Value = 'Los Lagos'
$ws.Range("D301").Value = 44358
$ws.Range("E301").Value = 10
$ws.Range("F301").Value = 'Fruta'
$ws.Range("G301").Value = 100104
$ws.Range("H301").Value = 'Frutos de pepita'
$ws.Range("I301").Value = 100104005
$ws.Range("J301").Value = 'Pera'
$ws.Range("K301").Value = "Packham's Triumph"
$ws.Range("L301").Value = 'Primera'
$ws.Range("M301").Value = 400
$ws.Range("N301").Value = 13500
$ws.Range("O301").Value = 14000
$ws.Range("P301").Value = 13750
$ws.Range("Q301").Value = '$/caja 15 kilos empedrada'
$ws.Range("R301").Value = "Región de O'Higgins"
$ws.Range("S301").Value = 917
$ws.Range("T301").Value = 15
